$wb = $excel.ActiveWorkbook

# Update the "zh-cn" worksheet: row 3 Correspond Handoff/Handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-11 13:38:57"
$wsZhCn.Range("G3").Value = "2016-01-11 13:40:02"

# Update the "de-de" worksheet: row 3 Correspond Handoff/Handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-11 13:39:15"
$wsDeDe.Range("G3").Value = "2016-01-11 13:40:33"
